$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: ZWO filter wheel part itself ---
# Part # (A) and Description (C) both become the new model name.
$ws.Range("A8").Value = "ZWO EFW 7x36mm"
$ws.Range("C8").Value = "ZWO EFW 7x36mm"
# Vendor (B8) and Comment (G8) stay the same ("ZWO" / "Check local distributors...").
# Unit price goes from 200 to 299 (F8 = E8*D8 recalculates automatically).
$ws.Range("E8").Value = 299

# --- Row 9: 3D printed filter adapter, resized for the new filter ---
$ws.Range("A9").Value = "32mmFilterAdapter.stl"
$ws.Range("G9").Value = "TODO for 32 mm filter in 36 mm slots (available 25mm filter for 31 mm slots). See /custom-parts/detection-accessories for printable 3D model"
$ws.Range("C9").Value = "Adapter to hold 32 mm filters in 36 mm slots"
# Vendor (B9) stays "3D printed (SLS)".

# --- New subtotal row for the "Detection: Tube lens and camera F-mount" section ---
$ws.Range("E32").Formula = "=SUM(E28:E31,E19,E9:E12)"

# Restore the selected cell shown in the saved workbook.
$ws.Range("C10").Select()
